# Auto-generated Excel COM-interop script
# Updates cached market-price / profit values across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled-runner scrape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1567.7142
$ws.Range("I38").Value = 1567.7142
$ws.Range("K38").Value = 4703.142599999999
$ws.Range("M38").Value = -4331.142599999999
$ws.Range("H64").Value = 4981.3438
$ws.Range("H67").Value = 4981.3438
$ws.Range("H69").Value = 26338.334
$ws.Range("I69").Value = 4500
$ws.Range("K69").Value = 13500
$ws.Range("M69").Value = -12626
$ws.Range("H72").Value = 26338.334
$ws.Range("I72").Value = 4500
$ws.Range("K72").Value = 40500
$ws.Range("M72").Value = -36132
$ws.Range("H87").Value = 119995
$ws.Range("J87").Value = 119995
$ws.Range("L87").Value = 119995
$ws.Range("N87").Value = -122491
$ws.Range("H90").Value = 119995
$ws.Range("J90").Value = 119995
$ws.Range("L90").Value = 359985
$ws.Range("N90").Value = -372465
$ws.Range("H96").Value = 1903.0714
$ws.Range("I96").Value = 1422.75
$ws.Range("K96").Value = 4268.25
$ws.Range("M96").Value = -2895.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1448.9333
$ws.Range("I2").Value = 1530.2727
$ws.Range("K2").Value = 1530.2727
$ws.Range("M2").Value = -1417.2727
$ws.Range("H74").Value = 6763540
$ws.Range("I74").Value = 10001660
$ws.Range("J74").Value = 17456.5
$ws.Range("K74").Value = 10001660
$ws.Range("L74").Value = 17456.5
$ws.Range("M74").Value = -10000786
$ws.Range("N74").Value = -19204.5
$ws.Range("H77").Value = 6763540
$ws.Range("I77").Value = 10001660
$ws.Range("J77").Value = 17456.5
$ws.Range("K77").Value = 50008300
$ws.Range("L77").Value = 87282.5
$ws.Range("M77").Value = -50003932
$ws.Range("N77").Value = -96018.5
$ws.Range("H97").Value = 1394.0416
$ws.Range("I97").Value = 971.7222
$ws.Range("J97").Value = 2661
$ws.Range("K97").Value = 971.7222
$ws.Range("L97").Value = 2661
$ws.Range("M97").Value = -475.7222
$ws.Range("N97").Value = -3653
$ws.Range("H116").Value = 1448.9333
$ws.Range("I116").Value = 1530.2727
$ws.Range("K116").Value = 1530.2727
$ws.Range("M116").Value = 763.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1448.9333
$ws.Range("I3").Value = 1530.2727
$ws.Range("K3").Value = 1530.2727
$ws.Range("M3").Value = -1416.2727
$ws.Range("H38").Value = 103000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H86").Value = 1390.9412
$ws.Range("I86").Value = 1319.5
$ws.Range("K86").Value = 1319.5
$ws.Range("M86").Value = -196.5
$ws.Range("H89").Value = 1390.9412
$ws.Range("I89").Value = 1319.5
$ws.Range("K89").Value = 6597.5
$ws.Range("M89").Value = -981.5
$ws.Range("H94").Value = 1957.5714
$ws.Range("I94").Value = 1375.75
$ws.Range("K94").Value = 1375.75
$ws.Range("M94").Value = -924.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2504.6
$ws.Range("I16").Value = 1839
$ws.Range("K16").Value = 1839
$ws.Range("M16").Value = -1552
$ws.Range("H99").Value = 4068.1333
$ws.Range("I99").Value = 3365.7273
$ws.Range("K99").Value = 3365.7273
$ws.Range("M99").Value = -1867.7273
$ws.Range("H113").Value = 2504.6
$ws.Range("I113").Value = 1839
$ws.Range("K113").Value = 1839
$ws.Range("M113").Value = 331
$ws.Range("H126").Value = 4068.1333
$ws.Range("I126").Value = 3365.7273
$ws.Range("K126").Value = 10097.1819
$ws.Range("M126").Value = -7627.1819
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 331.63635
$ws.Range("J23").Value = 260.42856
$ws.Range("L23").Value = 781.28568
$ws.Range("N23").Value = -1251.28568
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224
$ws.Range("H68").Value = 3057.375
$ws.Range("I68").Value = 3400
$ws.Range("K68").Value = 10200
$ws.Range("M68").Value = -9389
$ws.Range("H71").Value = 3057.375
$ws.Range("I71").Value = 3400
$ws.Range("K71").Value = 30600
$ws.Range("M71").Value = -26544
$ws.Range("H80").Value = 3641.5557
$ws.Range("J80").Value = 3641.5557
$ws.Range("L80").Value = 10924.6671
$ws.Range("N80").Value = -12796.6671
$ws.Range("H83").Value = 3641.5557
$ws.Range("J83").Value = 3641.5557
$ws.Range("L83").Value = 32774.0013
$ws.Range("N83").Value = -42134.0013
$ws.Range("H113").Value = 2523.1667
$ws.Range("I113").Value = 347.5
$ws.Range("K113").Value = 1042.5
$ws.Range("M113").Value = 1127.5
$ws.Range("H114").Value = 2109.5
$ws.Range("I114").Value = 2205.2
$ws.Range("J114").Value = 1631
$ws.Range("K114").Value = 6615.599999999999
$ws.Range("L114").Value = 4893
$ws.Range("M114").Value = -3361.599999999999
$ws.Range("N114").Value = -11401
$ws.Range("H122").Value = 1126.0526
$ws.Range("I122").Value = 697.6
$ws.Range("K122").Value = 6278.400000000001
$ws.Range("M122").Value = -3828.400000000001
$ws.Range("H132").Value = 2410.111
$ws.Range("I132").Value = 2139.9092
$ws.Range("J132").Value = 2834.7144
$ws.Range("K132").Value = 19259.1828
$ws.Range("L132").Value = 25512.4296
$ws.Range("M132").Value = -16729.1828
$ws.Range("N132").Value = -30572.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 48000
$ws.Range("I18").Value = 48000
$ws.Range("K18").Value = 48000
$ws.Range("M18").Value = -47707
$ws.Range("H126").Value = 4702.4375
$ws.Range("I126").Value = 4581.6665
$ws.Range("J126").Value = 4857.7144
$ws.Range("K126").Value = 13744.9995
$ws.Range("L126").Value = 14573.1432
$ws.Range("M126").Value = -11274.9995
$ws.Range("N126").Value = -19513.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 45455616
$ws.Range("I93").Value = 52632612
$ws.Range("K93").Value = 52632612
$ws.Range("M93").Value = -52631364
$ws.Range("H100").Value = 3387.8667
$ws.Range("I100").Value = 3618.1667
$ws.Range("K100").Value = 3618.1667
$ws.Range("M100").Value = -3077.1667
$ws.Range("H132").Value = 136124.94
$ws.Range("I132").Value = 114231.78
$ws.Range("J132").Value = 168964.67
$ws.Range("K132").Value = 342695.34
$ws.Range("L132").Value = 506894.01
$ws.Range("M132").Value = -340165.34
$ws.Range("N132").Value = -511954.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H81").Value = 50524
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 50524
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 17857870
$ws.Range("I107").Value = 27778664
$ws.Range("J107").Value = 444.5
$ws.Range("K107").Value = 83335992
$ws.Range("L107").Value = 1333.5
$ws.Range("M107").Value = -83334072
$ws.Range("N107").Value = -5173.5
$ws.Range("H136").Value = 1283
$ws.Range("I136").Value = 1283
$ws.Range("K136").Value = 3849
$ws.Range("M136").Value = -1299
